$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 9
$ws.Range("H9").Value = 263.66666
$ws.Range("J9").Value = 188.66667
$ws.Range("L9").Value = 188.66667
$ws.Range("N9").Value = -526.6666700000001
# Row 40
$ws.Range("H40").Value = 3600.3
$ws.Range("I40").Value = 1925.25
$ws.Range("J40").Value = 4717
$ws.Range("K40").Value = 1925.25
$ws.Range("L40").Value = 4717
$ws.Range("M40").Value = -1750.25
$ws.Range("N40").Value = -5067
# Row 55
$ws.Range("H55").Value = 299
$ws.Range("I55").Value = 400
$ws.Range("J55").Value = 198
$ws.Range("K55").Value = 400
$ws.Range("L55").Value = 198
$ws.Range("M55").Value = -186
$ws.Range("N55").Value = -626
# Row 137
$ws.Range("H137").Value = 2975.4443
$ws.Range("I137").Value = 2463.1667
$ws.Range("K137").Value = 7389.500100000001
$ws.Range("M137").Value = -4839.500100000001
# Row 138
$ws.Range("H138").Value = 2566.1853
$ws.Range("I138").Value = 2645.3076
$ws.Range("J138").Value = 2492.7144
$ws.Range("K138").Value = 7935.9228
$ws.Range("L138").Value = 7478.1432
$ws.Range("M138").Value = -2795.9228
$ws.Range("N138").Value = -17758.1432
# Row 141
$ws.Range("H141").Value = 5042.467
$ws.Range("I141").Value = 3617
$ws.Range("J141").Value = 24999
$ws.Range("K141").Value = 10851
$ws.Range("L141").Value = 74997
$ws.Range("M141").Value = -5671
$ws.Range("N141").Value = -85357

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 1858.5217
$ws.Range("I2").Value = 1907.4117
$ws.Range("J2").Value = 1720
$ws.Range("K2").Value = 1907.4117
$ws.Range("L2").Value = 1720
$ws.Range("M2").Value = -1794.4117
$ws.Range("N2").Value = -1946
# Row 32
$ws.Range("H32").Value = 2254.8733
$ws.Range("I32").Value = 2238.3286
$ws.Range("K32").Value = 2238.3286
$ws.Range("M32").Value = -1951.3286
# Row 63
$ws.Range("H63").Value = 3315
$ws.Range("I63").Value = 3315
$ws.Range("K63").Value = 3315
$ws.Range("M63").Value = -2629
# Row 66
$ws.Range("H66").Value = 3315
$ws.Range("I66").Value = 3315
$ws.Range("K66").Value = 16575
$ws.Range("M66").Value = -13143
# Row 81
$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = ""
$ws.Range("N81").Value = 0
# Row 84
$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = ""
$ws.Range("N84").Value = 0
# Row 97
$ws.Range("H97").Value = 1248.6666
$ws.Range("I97").Value = 1410.2354
$ws.Range("K97").Value = 1410.2354
$ws.Range("M97").Value = -914.2354
# Row 116
$ws.Range("H116").Value = 1858.5217
$ws.Range("I116").Value = 1907.4117
$ws.Range("J116").Value = 1720
$ws.Range("K116").Value = 1907.4117
$ws.Range("L116").Value = 1720
$ws.Range("M116").Value = 386.5882999999999
$ws.Range("N116").Value = -6308
# Row 132
$ws.Range("H132").Value = 2485.8823
$ws.Range("I132").Value = 2176.0715
$ws.Range("J132").Value = 3931.6667
$ws.Range("K132").Value = 6528.2145
$ws.Range("L132").Value = 11795.0001
$ws.Range("M132").Value = -3998.2145
$ws.Range("N132").Value = -16855.0001

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 1858.5217
$ws.Range("I3").Value = 1907.4117
$ws.Range("J3").Value = 1720
$ws.Range("K3").Value = 1907.4117
$ws.Range("L3").Value = 1720
$ws.Range("M3").Value = -1793.4117
$ws.Range("N3").Value = -1948
# Row 22
$ws.Range("H22").Value = 3732.4167
$ws.Range("I22").Value = 3732.4167
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 3732.4167
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = ""
$ws.Range("N22").Value = -3559.4167
# Row 86
$ws.Range("H86").Value = 4045.6365
$ws.Range("I86").Value = 4045.6365
$ws.Range("K86").Value = 4045.6365
$ws.Range("M86").Value = -2922.6365
# Row 89
$ws.Range("H89").Value = 4045.6365
$ws.Range("I89").Value = 4045.6365
$ws.Range("K89").Value = 20228.1825
$ws.Range("M89").Value = -14612.1825
# Row 94
$ws.Range("H94").Value = 9999
$ws.Range("I94").Value = 9999
$ws.Range("K94").Value = 9999
$ws.Range("M94").Value = -9548
# Row 105
$ws.Range("H105").Value = 7584.3335
$ws.Range("I105").Value = 11991.8
$ws.Range("J105").Value = 2075
$ws.Range("K105").Value = 11991.8
$ws.Range("L105").Value = 2075
$ws.Range("M105").Value = -10244.8
$ws.Range("N105").Value = -5569
# Row 134
$ws.Range("H134").Value = 131863.7
$ws.Range("I134").Value = 188509.5
$ws.Range("J134").Value = 2387.5715
$ws.Range("K134").Value = 565528.5
$ws.Range("L134").Value = 7162.7145
$ws.Range("M134").Value = -562993.5
$ws.Range("N134").Value = -12232.7145

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 54695.31
$ws.Range("I31").Value = 70120.53
$ws.Range("J31").Value = 33660.91
$ws.Range("K31").Value = 70120.53
$ws.Range("L31").Value = 33660.91
$ws.Range("M31").Value = -69825.53
$ws.Range("N31").Value = -34250.91
# Row 34
$ws.Range("H34").Value = 54695.31
$ws.Range("I34").Value = 70120.53
$ws.Range("J34").Value = 33660.91
$ws.Range("K34").Value = 70120.53
$ws.Range("L34").Value = 33660.91
$ws.Range("M34").Value = -69918.53
$ws.Range("N34").Value = -34064.91
# Row 62
$ws.Range("H62").Value = 55000
$ws.Range("I62").Value = 55000
$ws.Range("K62").Value = 55000
$ws.Range("M62").Value = -54376
# Row 65
$ws.Range("H65").Value = 55000
$ws.Range("I65").Value = 55000
$ws.Range("K65").Value = 275000
$ws.Range("M65").Value = -271880
# Row 134
$ws.Range("H134").Value = 3589.4211
$ws.Range("I134").Value = 3292.3333
$ws.Range("J134").Value = 4703.5
$ws.Range("K134").Value = 9876.999899999999
$ws.Range("L134").Value = 14110.5
$ws.Range("M134").Value = -7341.999899999999
$ws.Range("N134").Value = -19180.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 56
$ws.Range("H56").Value = 8253.143
$ws.Range("I56").Value = 8253.143
$ws.Range("K56").Value = 8253.143
$ws.Range("M56").Value = -7723.143
# Row 107
$ws.Range("H107").Value = 797.5
$ws.Range("J107").Value = 200
$ws.Range("L107").Value = 600
$ws.Range("N107").Value = -4440

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Range("H2").Value = 1218.375
$ws.Range("I2").Value = 587.5
$ws.Range("J2").Value = 1849.25
$ws.Range("K2").Value = 587.5
$ws.Range("L2").Value = 1849.25
$ws.Range("M2").Value = -474.5
$ws.Range("N2").Value = -2075.25
# Row 70
$ws.Range("H70").Value = 2428
$ws.Range("I70").Value = 2428
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 2428
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = ""
$ws.Range("N70").Value = -2158
# Row 73
$ws.Range("H73").Value = 2428
$ws.Range("I73").Value = 2428
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 2428
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = ""
$ws.Range("N73").Value = -1492
# Row 97
$ws.Range("H97").Value = 88830.47
$ws.Range("I97").Value = 61474.91
$ws.Range("K97").Value = 61474.91
$ws.Range("M97").Value = -60978.91

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 4700
$ws.Range("I22").Value = 1750
$ws.Range("J22").Value = 6666.6665
$ws.Range("K22").Value = 1750
$ws.Range("L22").Value = 6666.6665
$ws.Range("M22").Value = -1455
$ws.Range("N22").Value = -7256.6665
# Row 27
$ws.Range("H27").Value = 4700
$ws.Range("I27").Value = 1750
$ws.Range("J27").Value = 6666.6665
$ws.Range("K27").Value = 1750
$ws.Range("L27").Value = 6666.6665
$ws.Range("M27").Value = -1643
$ws.Range("N27").Value = -6880.6665
# Row 42
$ws.Range("H42").Value = 18999.5
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 18999.5
$ws.Range("K42").Value = 0
$ws.Range("L42").Value = ""
$ws.Range("M42").Value = 18999.5
$ws.Range("N42").Value = -20125.5
# Row 46
$ws.Range("H46").Value = 3637.3333
$ws.Range("I46").Value = 3306.8
$ws.Range("J46").Value = 4050.5
$ws.Range("K46").Value = 3306.8
$ws.Range("L46").Value = 4050.5
$ws.Range("M46").Value = -3118.8
$ws.Range("N46").Value = -4426.5
# Row 49
$ws.Range("H49").Value = 18999.5
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 18999.5
$ws.Range("K49").Value = 0
$ws.Range("L49").Value = ""
$ws.Range("M49").Value = 18999.5
$ws.Range("N49").Value = -19293.5
# Row 122
$ws.Range("H122").Value = 129125.81
$ws.Range("I122").Value = 157462.53
$ws.Range("K122").Value = 472387.59
$ws.Range("M122").Value = -469937.59
# Row 132
$ws.Range("H132").Value = 3097.8235
$ws.Range("I132").Value = 2727.0344
$ws.Range("J132").Value = 5248.4
$ws.Range("K132").Value = 8181.1032
$ws.Range("L132").Value = 15745.2
$ws.Range("M132").Value = -5651.1032
$ws.Range("N132").Value = -20805.2
# Row 136
$ws.Range("H136").Value = 5071.8
$ws.Range("I136").Value = 4706.4
$ws.Range("J136").Value = 6898.8
$ws.Range("K136").Value = 14119.2
$ws.Range("L136").Value = 20696.4
$ws.Range("M136").Value = -11569.2
$ws.Range("N136").Value = -25796.4

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 70
$ws.Range("H70").Value = 20089.285
$ws.Range("J70").Value = 20089.285
$ws.Range("L70").Value = 20089.285
$ws.Range("N70").Value = -20719.285
# Row 73
$ws.Range("H73").Value = 20089.285
$ws.Range("J73").Value = 20089.285
$ws.Range("L73").Value = 20089.285
$ws.Range("N73").Value = -22273.285
# Row 81
$ws.Range("H81").Value = 4029.2
$ws.Range("I81").Value = 950.5
$ws.Range("K81").Value = 1901
$ws.Range("M81").Value = -840
# Row 84
$ws.Range("H84").Value = 4029.2
$ws.Range("I84").Value = 950.5
$ws.Range("K84").Value = 9505
$ws.Range("M84").Value = -4201
# Row 100
$ws.Range("H100").Value = 1364.1052
$ws.Range("I100").Value = 965.26666
$ws.Range("J100").Value = 2859.75
$ws.Range("K100").Value = 1930.53332
$ws.Range("L100").Value = 5719.5
$ws.Range("M100").Value = -1389.53332
$ws.Range("N100").Value = -6801.5
# Row 132
$ws.Range("H132").Value = 8480.182000000001
$ws.Range("I132").Value = 11833.214
$ws.Range("J132").Value = 2612.375
$ws.Range("K132").Value = 35499.642
$ws.Range("L132").Value = 7837.125
$ws.Range("M132").Value = -32969.642
$ws.Range("N132").Value = -12897.125
